$d = $word.ActiveDocument

# 1. "Course Title: Learn MongoDB" -- merge the trailing " " and "Learn MongoDB"
#    runs into a single run " Learn MongoDB".
$d.Content.Find.Execute("Learn MongoDB", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Learn MongoDB", 2) | Out-Null

# 2. "Provider: codeacademy.com" -- merge the trailing " " and "codeacademy.com"
#    runs into a single run " codeacademy.com".
$d.Content.Find.Execute("codeacademy.com", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "codeacademy.com", 2) | Out-Null

# 3. "Progress: 26%" -> "Progress: 52%", splitting into ": ", "52", "%" runs
#    (mirrors someone retyping just the "26" -> "52" digits in place).
$progressRange = $d.Content.Duplicate
$progressRange.Find.Execute("26", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0) | Out-Null
$progressRange.Text = "52"
# Toggling Bold on/off on the freshly-typed range forces Word to keep it as
# its own run (instead of silently re-coalescing with neighboring runs that
# share identical formatting), without leaving any stray formatting behind.
$progressRange.Bold = 1
$progressRange.Bold = 0

# 4. "Duration: 5 hours" -- merge the " ", "5", " ", "hours" runs into a
#    single run " 5 hours".
$d.Content.Find.Execute("5 hours", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "5 hours", 2) | Out-Null

# 5. Add a new bullet after the ".find()" bullet describing .insertOne.
$findRange = $d.Content.Duplicate
$findRange.Find.Execute("method retrieves documents from a collection and returns a cursor", `
                         $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$findRange.Collapse(0)
$findRange.InsertParagraphAfter()
$findRange.Collapse(0)
$findRange.Move(4, 1) | Out-Null
$findRange.Text = "Utilized the .insertOne method to efficiently insert a single document into a collection"
